$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.999.95"
$ws.Range("E2").Value = "  -0.46%  "
$ws.Range("D3").Value = "3.473.11"
$ws.Range("E3").Value = "  -1.48%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "592.76"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.69%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "176.19"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.60%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.587"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.56%  "
$ws.Range("E9").Value = "  -2.77%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.08"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.10%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.426"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.94%  "
$ws.Range("D12").Value = "4.071.10"
$ws.Range("E12").Value = "  -1.59%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "30.83"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +7.46%  "
$ws.Range("E14").Value = "  -0.11%  "
$ws.Range("D15").Value = "67.096.86"
$ws.Range("E15").Value = "  -0.18%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000176"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.98%  "
$ws.Range("D17").Value = "3.460.11"
$ws.Range("E17").Value = "  -1.93%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.24"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.77%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.33"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.24%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "387.99"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.01%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.87"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.34%  "
$ws.Range("E22").Value = "  +0.19%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "72.60"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.16%  "
$ws.Range("E24").Value = "  -0.22%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.534"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.08%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000121"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.90%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.29"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.53%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.178"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.28%  "
$ws.Range("E29").Value = "  -0.48%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.10"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.12%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.42"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.32%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.04"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.21%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "23.45"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.55%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.27"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.50%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.62"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.22%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "163.20"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.47%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.873"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.90%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.92"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.34%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.92"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.31%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "27.33"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.43%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.62"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.19%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "26.29"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.86%  "
$ws.Range("D43").Value = "2.789.50"
$ws.Range("E43").Value = "  -0.81%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0722"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.32%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.57"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.65%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "42.18"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.83%  "
$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0299"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.41%  "
$ws.Range("B48").Value = "Bittensor"
$ws.Range("C48").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "339.17"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.92%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.07"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.32%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "33.11"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.34%  "
$ws.Range("B51").Value = "Cosmos"
$ws.Range("C51").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.37"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.39%  "
